# Auto-generated script applying the GitHub Actions cryptos update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.732.06"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "1.892.35"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'244.86"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4921"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.2960"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'0.06794"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").Value = "1.890.36"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("D12").Value = "'0.07236"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "'91.26"
$ws.Range("E13").Value = "  +5.62%  "
$ws.Range("D14").Value = "'0.6805"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "'5.051"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "30.676.20"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "'0.000008008"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("D20").Value = "2.131.80"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'4.825"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "'190.72"
$ws.Range("E23").Value = "  +33.09%  "
$ws.Range("D24").Value = "'6.129"
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("D25").Value = "'9.371"
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").Value = "'154.74"
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D27").Value = "'19.14"
$ws.Range("E27").Value = "  +12.95%  "
$ws.Range("D28").Value = "'1.905"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "'1.402"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Value = "'4.351"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("D31").Value = "'0.09096"
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("D32").Value = "'4.018"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "'0.05205"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("D34").Value = "'0.7517"
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +4.11%  "
$ws.Range("D37").Value = "'0.01846"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "'0.9396"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'0.4434"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("D42").Value = "'105.39"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "'5.775"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D45").Value = "'7.617"
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").Value = "'0.1347"
$ws.Range("E46").Value = "  +6.13%  "
$ws.Range("D47").Value = "'0.05867"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").Value = "'8.756"
$ws.Range("E48").Value = "  +5.93%  "
$ws.Range("E49").Value = "  +6.35%  "
$ws.Range("D50").Value = "'0.3943"
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").Value = "'33.67"
$ws.Range("E51").Value = "  +2.84%  "
